$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Baseline-Config")

# Update the historical end year from 2018 to 2019 (DOLA 2020 dataset).
# Dependent cells (E17, F16, F17) recalc automatically via formulas.
$ws.Range("E16").Value = 2019

# Apply the "shaded" style (same as used on E15/E17/E18, style index 12 -
# light gray fill) to the formula-text column F10:F18.
$ws.Range("E15").Copy() | Out-Null
$ws.Range("F10:F18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selected/active cell shown when the sheet is reopened.
$ws.Activate()
$ws.Range("E25").Select() | Out-Null
